$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has formatted (but empty) placeholder rows at 18 and 19.
# Copy the formatting from the last "normal" data row (16) down onto rows
# 18:19 so the new entries look consistent with the rest of the table,
# then fill in the new log entries.
$ws.Range("A16:E16").Copy()
$ws.Range("A18:E19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 18 - Sprites and Animations
$ws.Range("A18").Value = "Sprites and Animations"
$ws.Range("B18").Value = 44323.0
$ws.Range("C18").Value = 0.0375
$ws.Range("D18").Value = 0.07361111111111111
$ws.Range("E18").Value = 0.036006944444444446

# Row 19 - Settings Window & Collision Fixes
$ws.Range("A19").Value = "Settings Window & Collision Fixes"
$ws.Range("B19").Value = 44323.0
$ws.Range("C19").Value = 0.49722222222222223
$ws.Range("D19").Value = 0.53125
$ws.Range("E19").Value = 0.03350694444444444
